$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7965988516807556
$ws.Range("B1").Value = 1.16258704662323
$ws.Range("C1").Value = 3.5346839427948
$ws.Range("D1").Value = 3.872872114181519
$ws.Range("E1").Value = 1.627381801605225
